# Revert "reapply some glitch changes":
#   - Span Example (slide 3) + Div Example (slide 5): swap the glitch.com
#     remix links back to the replit.com links, dropping the stray
#     trailing-space run that used to follow the hyperlink run.
#   - slide layout 8's cached "today" date field text.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide layout 8 (used by slide 3) - cached datetimeFigureOut text.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$dateShape = $s3.CustomLayout.Shapes.Item(3)   # "Date Placeholder 3"
$dateShape.TextFrame.TextRange.Text = "2/9/2023"

# ---------------------------------------------------------------------
# 2) Slide 3 "Span Example" - Content Placeholder 2.
#    Run 1 (chars 1-44) = the hyperlinked URL text, keep its rPr
#    (hyperlink/highlight/bold/size/color) but change the visible text.
#    Run 2 (the trailing " ") gets deleted entirely.
# ---------------------------------------------------------------------
$spanShape = $s3.Shapes.Item(2)   # "Content Placeholder 2"
$spanRange = $spanShape.TextFrame.TextRange

$spanUrlLen = "https://glitch.com/edit/#!/remix/spanexample".Length
$spanUrlRun = $spanRange.Characters(1, $spanUrlLen)
$spanUrlRun.Text = "https://replit.com/@HylandOutreach/SpanExample"

$spanNewLen = $spanRange.Length
$spanTailRun = $spanRange.Characters($spanNewLen, 1)
$spanTailRun.Text = ""

# ---------------------------------------------------------------------
# 3) Slide 5 "Div Example" - Rectangle 3 shape.
#    Same pattern: keep run 1's rPr, change its text; drop the trailing
#    " " run.
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$divShape = $s5.Shapes.Item(3)   # "Rectangle 3"
$divRange = $divShape.TextFrame.TextRange

$divUrlLen = "https://glitch.com/edit/#!/remix/divexample".Length
$divUrlRun = $divRange.Characters(1, $divUrlLen)
$divUrlRun.Text = "https://replit.com/@HylandOutreach/DivExample"

$divNewLen = $divRange.Length
$divTailRun = $divRange.Characters($divNewLen, 1)
$divTailRun.Text = ""
